$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers (U = "3-jul", V = "4-jul")
$ws.Range("U1").Value = "3-jul"
$ws.Range("V1").Value = "4-jul"

# Data for the two new day columns, rows 2-18
$uValues = @(0, 15, 14, 30, 0, 29, 10, 22, 28, 14, 0, 13, 0, 0, 17, 0, 0)
$vValues = @(0, 15.013077939530611, 14.196482638689593, 29.818282710578451, 0, 24.394734874761802, 14.931930065813283, 21.771032694534203, 27.669167259839902, 13.812642495354032, 0, 14.455591408120702, 0, 0, 15.920611459455261, 0, 0)

for ($i = 0; $i -lt $uValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 21).Value = $uValues[$i]
    $ws.Cells.Item($row, 21).HorizontalAlignment = -4108
    $ws.Cells.Item($row, 22).Value = $vValues[$i]
}

# Match the selection recorded in the saved view
$ws.Range("P11:P12").Select()
